$d = $word.ActiveDocument

# 1. Add "Miñano, Richard" as a new list item before "Suarez, Santiago"
$rngSuarez = $d.Content
$rngSuarez.Find.Execute("Suarez, Santiago")
$paraSuarez = $rngSuarez.Paragraphs(1)
$paraSuarez.Range.InsertParagraphBefore()

$rngNew = $d.Content
$rngNew.Find.Execute("Suarez, Santiago")
$paraNew = $rngNew.Paragraphs(1).Previous()
$paraNew.Range.InsertBefore("Miñano, Richard")

# 2. Merge the "...uso." paragraph with the following page-break-only paragraph
#    (delete the paragraph mark between them so the page break run joins the
#    previous paragraph).
$rngUso = $d.Content
$rngUso.Find.Execute("fácil comprensión y uso.")
$paraUso = $rngUso.Paragraphs(1)
$mark = $d.Range($paraUso.Range.End - 1, $paraUso.Range.End)
$mark.Delete()

# 3. Wrap "Trello" with a _GoBack bookmark, splitting its run in three
$rngTrello = $d.Content
$rngTrello.Find.Execute("Trello")
$d.Bookmarks.Add("_GoBack", $rngTrello)

# 4. The old _GoBack bookmark (previously sitting between ")" and
#    " que lo respalden.") was already relocated by step 3's Bookmarks.Add.
#    Now merge ")" with " que lo respalden." back into a single run.
$rngTail = $d.Content
$rngTail.Find.Execute("etc.) que lo respalden.")
$parenStart = $rngTail.Start + 4

$tailRng = $d.Range($parenStart, $rngTail.End)
$tailRng.Delete()
$insertPoint = $d.Range($parenStart, $parenStart)
$insertPoint.InsertAfter(") que lo respalden.")
